$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0736"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.60"
$ws.Range("E41").Value = "  +5.54%  "

$ws.Range("D2").Value = "66.428.70"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").Value = "3.488.41"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.59"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.58"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.483.66"
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  +6.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "4.098.33"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.92"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "66.556.05"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "3.489.91"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").Value = "  +3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.72"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.85"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.71"
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.529"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  +4.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +6.49%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.56"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.29"
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +6.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.98"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.892"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.78"
$ws.Range("E39").Value = "  +5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.26"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "2.776.22"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.40"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.73"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0308"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "341.15"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.26"
$ws.Range("E50").Value = "  +10.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.851"
$ws.Range("E51").Value = "  +4.61%  "
